$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A2")

# Force a text entry so the new SKU's leading zero is preserved, then
# strip back down to the sheet's plain/default style (borrowed from a
# neighboring cell that already carries no special number format) so the
# cell ends up with no explicit "text" formatting applied to it.
$cell.NumberFormat = "@"
$cell.Value = "0191448237469"
$cell.Style = $ws.Range("B2").Style
